$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing F/G values per diff
$ws.Cells.Item(306, 6).Value = 74521
$ws.Cells.Item(344, 6).Value = 136303
$ws.Cells.Item(344, 7).Value = 2475
$ws.Cells.Item(347, 6).Value = 345903
$ws.Cells.Item(371, 6).Value = 160475
$ws.Cells.Item(371, 7).Value = 1976
$ws.Cells.Item(393, 6).Value = 309627
$ws.Cells.Item(400, 6).Value = 150086
$ws.Cells.Item(414, 6).Value = 149324
$ws.Cells.Item(428, 6).Value = 102746
$ws.Cells.Item(454, 6).Value = 52785
$ws.Cells.Item(455, 6).Value = 50891
$ws.Cells.Item(462, 6).Value = 43806
$ws.Cells.Item(667, 6).Value = 16847
$ws.Cells.Item(677, 6).Value = 56401
$ws.Cells.Item(677, 7).Value = 802
$ws.Cells.Item(687, 6).Value = 31505
$ws.Cells.Item(692, 6).Value = 41691
$ws.Cells.Item(705, 6).Value = 56021
$ws.Cells.Item(706, 6).Value = 40678
$ws.Cells.Item(708, 6).Value = 35603
$ws.Cells.Item(709, 6).Value = 32387
$ws.Cells.Item(712, 6).Value = 51389
$ws.Cells.Item(714, 6).Value = 32519
$ws.Cells.Item(722, 6).Value = 27947
$ws.Cells.Item(723, 6).Value = 22561
$ws.Cells.Item(727, 6).Value = 25099
$ws.Cells.Item(728, 6).Value = 24661
$ws.Cells.Item(729, 6).Value = 23228
$ws.Cells.Item(730, 6).Value = 19472
$ws.Cells.Item(730, 7).Value = 2328
$ws.Cells.Item(733, 6).Value = 31637
$ws.Cells.Item(733, 7).Value = 3715
$ws.Cells.Item(734, 6).Value = 23138
$ws.Cells.Item(734, 7).Value = 2547
$ws.Cells.Item(735, 6).Value = 19307
$ws.Cells.Item(735, 7).Value = 2273
$ws.Cells.Item(736, 6).Value = 19533
$ws.Cells.Item(736, 7).Value = 2186
$ws.Cells.Item(737, 6).Value = 18457
$ws.Cells.Item(739, 6).Value = 8645
$ws.Cells.Item(739, 7).Value = 1398
$ws.Cells.Item(740, 6).Value = 24574
$ws.Cells.Item(740, 7).Value = 2735
$ws.Cells.Item(741, 6).Value = 18851
$ws.Cells.Item(741, 7).Value = 1915
$ws.Cells.Item(742, 6).Value = 17223
$ws.Cells.Item(742, 7).Value = 1677
$ws.Cells.Item(743, 6).Value = 17958
$ws.Cells.Item(743, 7).Value = 1603
$ws.Cells.Item(744, 6).Value = 14655
$ws.Cells.Item(744, 7).Value = 1597
$ws.Cells.Item(745, 6).Value = 6119
$ws.Cells.Item(745, 7).Value = 914
$ws.Cells.Item(746, 6).Value = 7888
$ws.Cells.Item(746, 7).Value = 1223
$ws.Cells.Item(747, 6).Value = 22013
$ws.Cells.Item(747, 7).Value = 2322
$ws.Cells.Item(748, 6).Value = 16756
$ws.Cells.Item(748, 7).Value = 1508
$ws.Cells.Item(749, 6).Value = 14638
$ws.Cells.Item(749, 7).Value = 1450
$ws.Cells.Item(750, 6).Value = 14837
$ws.Cells.Item(750, 7).Value = 1316
$ws.Cells.Item(751, 6).Value = 12405
$ws.Cells.Item(751, 7).Value = 1352
$ws.Cells.Item(752, 6).Value = 4717
$ws.Cells.Item(752, 7).Value = 605
$ws.Cells.Item(753, 6).Value = 6431
$ws.Cells.Item(753, 7).Value = 912
$ws.Cells.Item(754, 6).Value = 20700
$ws.Cells.Item(754, 7).Value = 1889
$ws.Cells.Item(755, 6).Value = 13502
$ws.Cells.Item(755, 7).Value = 1268

# Add new rows 756-760
$ws.Cells.Item(756, 1).Value = 44650
$ws.Cells.Item(756, 1).NumberFormat = "yyyy-mm-dd"
$ws.Cells.Item(756, 2).Value = 1710884
$ws.Cells.Item(756, 3).Value = 15828
$ws.Cells.Item(756, 4).Value = 7334
$ws.Cells.Item(756, 5).Value = 19352
$ws.Cells.Item(756, 6).Value = 13394
$ws.Cells.Item(756, 7).Value = 1030

$ws.Cells.Item(757, 1).Value = 44651
$ws.Cells.Item(757, 1).NumberFormat = "yyyy-mm-dd"
$ws.Cells.Item(757, 2).Value = 1716638
$ws.Cells.Item(757, 3).Value = 11822
$ws.Cells.Item(757, 4).Value = 5754
$ws.Cells.Item(757, 5).Value = 19368
$ws.Cells.Item(757, 6).Value = 12856
$ws.Cells.Item(757, 7).Value = 949

$ws.Cells.Item(758, 1).Value = 44652
$ws.Cells.Item(758, 1).NumberFormat = "yyyy-mm-dd"
$ws.Cells.Item(758, 2).Value = 1722229
$ws.Cells.Item(758, 3).Value = 13002
$ws.Cells.Item(758, 4).Value = 5591
$ws.Cells.Item(758, 5).Value = 19392
$ws.Cells.Item(758, 6).Value = 9073
$ws.Cells.Item(758, 7).Value = 818

$ws.Cells.Item(759, 1).Value = 44653
$ws.Cells.Item(759, 1).NumberFormat = "yyyy-mm-dd"
$ws.Cells.Item(759, 2).Value = 1725487
$ws.Cells.Item(759, 3).Value = 8024
$ws.Cells.Item(759, 4).Value = 3258
$ws.Cells.Item(759, 5).Value = 19417
$ws.Cells.Item(759, 6).Value = 2814
$ws.Cells.Item(759, 7).Value = 312

$ws.Cells.Item(760, 1).Value = 44654
$ws.Cells.Item(760, 1).NumberFormat = "yyyy-mm-dd"
$ws.Cells.Item(760, 2).Value = 1726952
$ws.Cells.Item(760, 3).Value = 4041
$ws.Cells.Item(760, 4).Value = 1465
$ws.Cells.Item(760, 5).Value = 19440
$ws.Cells.Item(760, 6).Value = 3070
$ws.Cells.Item(760, 7).Value = 400

